# edit.ps1 - applies the target changes to the document via Word COM interop.
# Strategy: for each paragraph that needs modification, locate it with a unique
# Find query (using wildcards to skip over special characters such as en-dashes
# that don't always round-trip cleanly through the scripting host), then replace
# the *entire* paragraph's XML with InsertXML so the resulting run/proofErr
# structure exactly matches the target OOXML.

$d = $word.ActiveDocument

# --- Change 1: "Data Analyst Intern" date range: split " – " into " –" + " " ---
$rng1 = $d.Content
$rng1.Find.Execute("Jun. 2023*Present", $true, $false, $true, $false, $false, $true, 1, $false, $null, 0) | Out-Null
$xml1 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="0000002B" w14:textId="72E9344A" w:rsidR="003830FF" w:rsidRDefault="00000000"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="right"/><w:rPr><w:rFonts w:ascii="Georgia" w:eastAsia="Georgia" w:hAnsi="Georgia" w:cs="Georgia"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:sdt><w:sdtPr><w:tag w:val="goog_rdk_6"/><w:id w:val="855929888"/><w:showingPlcHdr/></w:sdtPr><w:sdtContent><w:r w:rsidR="00770715"><w:t xml:space="preserve">     </w:t></w:r></w:sdtContent></w:sdt><w:r w:rsidR="00770715"><w:t xml:space="preserve">        </w:t></w:r><w:r w:rsidR="00CF24F1"><w:rPr><w:rFonts w:ascii="Georgia" w:eastAsia="Georgia" w:hAnsi="Georgia" w:cs="Georgia"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Jun</w:t></w:r><w:r w:rsidR="002C385D"><w:rPr><w:rFonts w:ascii="Georgia" w:eastAsia="Georgia" w:hAnsi="Georgia" w:cs="Georgia"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>.</w:t></w:r><w:r w:rsidR="00CF24F1"><w:rPr><w:rFonts w:ascii="Georgia" w:eastAsia="Georgia" w:hAnsi="Georgia" w:cs="Georgia"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00CF24F1"><w:t>2</w:t></w:r><w:r w:rsidR="00CF24F1"><w:rPr><w:rFonts w:ascii="Georgia" w:eastAsia="Georgia" w:hAnsi="Georgia" w:cs="Georgia"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>023</w:t></w:r><w:r w:rsidR="002C385D"><w:rPr><w:rFonts w:ascii="Georgia" w:eastAsia="Georgia" w:hAnsi="Georgia" w:cs="Georgia"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> –</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Georgia" w:eastAsia="Georgia" w:hAnsi="Georgia" w:cs="Georgia"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00CF24F1"><w:rPr><w:rFonts w:ascii="Georgia" w:eastAsia="Georgia" w:hAnsi="Georgia" w:cs="Georgia"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Present</w:t></w:r><w:r w:rsidR="002C385D"><w:rPr><w:rFonts w:ascii="Georgia" w:eastAsia="Georgia" w:hAnsi="Georgia" w:cs="Georgia"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:sdt><w:sdtPr><w:tag w:val="goog_rdk_7"/><w:id w:val="-1483074941"/><w:showingPlcHdr/></w:sdtPr><w:sdtContent><w:r w:rsidR="00CF24F1"><w:t xml:space="preserve">     </w:t></w:r></w:sdtContent></w:sdt></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$rng1.InsertXML($xml1)

# --- Change 2: "AI Community Surfer" date range: "Present" -> "Jul 2023 " ---
$rng2 = $d.Content
$rng2.Find.Execute("May. 2023*Present", $true, $false, $true, $false, $false, $true, 1, $false, $null, 0) | Out-Null
$xml2 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="00000031" w14:textId="4796B281" w:rsidR="003830FF" w:rsidRDefault="00C9553F"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="right"/><w:rPr><w:rFonts w:ascii="Georgia" w:eastAsia="Georgia" w:hAnsi="Georgia" w:cs="Georgia"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Georgia" w:eastAsia="Georgia" w:hAnsi="Georgia" w:cs="Georgia"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">  </w:t></w:r><w:r w:rsidR="00770715"><w:rPr><w:rFonts w:ascii="Georgia" w:eastAsia="Georgia" w:hAnsi="Georgia" w:cs="Georgia"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">          </w:t></w:r><w:r w:rsidR="00CF24F1"><w:rPr><w:rFonts w:ascii="Georgia" w:eastAsia="Georgia" w:hAnsi="Georgia" w:cs="Georgia"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>May</w:t></w:r><w:r w:rsidR="002C385D"><w:rPr><w:rFonts w:ascii="Georgia" w:eastAsia="Georgia" w:hAnsi="Georgia" w:cs="Georgia"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>.</w:t></w:r><w:r w:rsidR="00CF24F1"><w:rPr><w:rFonts w:ascii="Georgia" w:eastAsia="Georgia" w:hAnsi="Georgia" w:cs="Georgia"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> 2023</w:t></w:r><w:r w:rsidR="002C385D"><w:rPr><w:rFonts w:ascii="Georgia" w:eastAsia="Georgia" w:hAnsi="Georgia" w:cs="Georgia"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> –</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Georgia" w:eastAsia="Georgia" w:hAnsi="Georgia" w:cs="Georgia"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Georgia" w:eastAsia="Georgia" w:hAnsi="Georgia" w:cs="Georgia"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Jul 202</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Georgia" w:eastAsia="Georgia" w:hAnsi="Georgia" w:cs="Georgia"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>3</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Georgia" w:eastAsia="Georgia" w:hAnsi="Georgia" w:cs="Georgia"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$rng2.InsertXML($xml2)

# --- Change 3: merge "compared" / gramStart-End runs into one run ---
$rng3 = $d.Content
$rng3.Find.Execute("compared", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0) | Out-Null
$xml3 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="27C2E35D" w14:textId="72626213" w:rsidR="003830FF" w:rsidRPr="00E05814" w:rsidRDefault="00E05814" w:rsidP="00E05814"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:pBdr><w:top w:val="nil"/><w:left w:val="nil"/><w:bottom w:val="nil"/><w:right w:val="nil"/><w:between w:val="nil"/></w:pBdr><w:tabs><w:tab w:val="num" w:pos="720"/></w:tabs><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Georgia" w:eastAsia="Georgia" w:hAnsi="Georgia" w:cs="Georgia"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Georgia" w:eastAsia="Georgia" w:hAnsi="Georgia" w:cs="Georgia"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">I trained 6 </w:t></w:r><w:r w:rsidRPr="00E05814"><w:rPr><w:rFonts w:ascii="Georgia" w:eastAsia="Georgia" w:hAnsi="Georgia" w:cs="Georgia"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>different machine learning algorithms</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Georgia" w:eastAsia="Georgia" w:hAnsi="Georgia" w:cs="Georgia"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>,</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Georgia" w:eastAsia="Georgia" w:hAnsi="Georgia" w:cs="Georgia"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> compared and evaluated </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Georgia" w:eastAsia="Georgia" w:hAnsi="Georgia" w:cs="Georgia"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">them </w:t></w:r><w:r w:rsidRPr="00E05814"><w:rPr><w:rFonts w:ascii="Georgia" w:eastAsia="Georgia" w:hAnsi="Georgia" w:cs="Georgia"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>using default hyperparameters.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Georgia" w:eastAsia="Georgia" w:hAnsi="Georgia" w:cs="Georgia"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> After which I fine-tuned the best performing model using </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Georgia" w:eastAsia="Georgia" w:hAnsi="Georgia" w:cs="Georgia"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>GridsearchCV</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Georgia" w:eastAsia="Georgia" w:hAnsi="Georgia" w:cs="Georgia"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> which gave it 2% more accuracy.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$rng3.InsertXML($xml3)

# --- Change 4: merge "dataset," / "and" gramStart-End runs into one run ---
$rng4 = $d.Content
$rng4.Find.Execute("I collaborated with a team", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0) | Out-Null
$xml4 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="43597727" w14:textId="51A35E37" w:rsidR="00E05814" w:rsidRPr="00E05814" w:rsidRDefault="00E05814" w:rsidP="00E05814"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:pBdr><w:top w:val="nil"/><w:left w:val="nil"/><w:bottom w:val="nil"/><w:right w:val="nil"/><w:between w:val="nil"/></w:pBdr><w:tabs><w:tab w:val="num" w:pos="720"/></w:tabs><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Georgia" w:eastAsia="Georgia" w:hAnsi="Georgia" w:cs="Georgia"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Georgia" w:eastAsia="Georgia" w:hAnsi="Georgia" w:cs="Georgia"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">I collaborated with a team to conduct an exploratory data analysis on the dataset, </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Georgia" w:eastAsia="Georgia" w:hAnsi="Georgia" w:cs="Georgia"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>and provided valuable insights and recommendations for improving customer retention in the telecommunications company.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$rng4.InsertXML($xml4)

Write-Host "All edits applied."
